$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row of data (row 10)
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Rachel"
$ws.Range("C10").Value = "Milo"

# Update the selected cell to match the new state
$ws.Range("H12").Select()
